# Natmi following Dr Hou advice
#
# The Edil3-Itgb3 sending/target cluster table is recomputed: the three
# sending clusters (FAPs, M2, sCs) are now each crossed with all four target
# clusters (ECs, FAPs, M2, sCs) - 12 data rows total (rows 2-13) instead of
# the previous 8 (rows 2-9, which only covered FAPs/sCs as senders).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A..T for each data row, in row order (row 2 .. row 13).
$rows = @(
    @("FAPs", "Edil3", "Itgb3", "ECs", 3, 1, 8.032126666666667, 24.09638, 0.8752332474713129, 0.8752332474713128, 3, 1, 3.778439, 11.335317, 0.4252971528324392, 0.4252971528324392, 30.34890065027333, 273.14010585246, 0.3722342082138391, 0.3722342082138391),
    @("FAPs", "Edil3", "Itgb3", "FAPs", 3, 1, 8.032126666666667, 24.09638, 0.8752332474713129, 0.8752332474713128, 3, 1, 4.333403333333333, 13.00021, 0.4877633593505858, 0.4877633593505858, 34.80644447108889, 313.2580002398, 0.4269067090019302, 0.4269067090019302),
    @("FAPs", "Edil3", "Itgb3", "M2", 3, 1, 8.032126666666667, 24.09638, 0.8752332474713129, 0.8752332474713128, 3, 1, 0.2909853333333334, 0.8729560000000001, 0.03275300561492853, 0.03275300561492853, 2.337231055475556, 21.03507949928, 0.02866651946880005, 0.02866651946880004),
    @("FAPs", "Edil3", "Itgb3", "sCs", 3, 1, 8.032126666666667, 24.09638, 0.8752332474713129, 0.8752332474713128, 3, 1, 0.4814053333333333, 1.444216, 0.0541864822020464, 0.05418648220204641, 3.866708615342222, 34.80037753808, 0.04742581078674357, 0.04742581078674357),
    @("M2", "Edil3", "Itgb3", "ECs", 1, 0.3333333333333333, 0.05236733333333333, 0.157102, 0.005706288398682217, 0.005706288398682216, 3, 1, 3.778439, 11.335317, 0.4252971528324392, 0.4252971528324392, 0.1978667745926667, 1.780800971334, 0.002426868209200326, 0.002426868209200325),
    @("M2", "Edil3", "Itgb3", "FAPs", 1, 0.3333333333333333, 0.05236733333333333, 0.157102, 0.005706288398682217, 0.005706288398682216, 3, 1, 4.333403333333333, 13.00021, 0.4877633593505858, 0.4877633593505858, 0.2269287768244444, 2.04235899142, 0.002783318398764513, 0.002783318398764513),
    @("M2", "Edil3", "Itgb3", "M2", 1, 0.3333333333333333, 0.05236733333333333, 0.157102, 0.005706288398682217, 0.005706288398682216, 3, 1, 0.2909853333333334, 0.8729560000000001, 0.03275300561492853, 0.03275300561492853, 0.01523812594577778, 0.137143133512, 0.0001868980959624402, 0.0001868980959624402),
    @("M2", "Edil3", "Itgb3", "sCs", 1, 0.3333333333333333, 0.05236733333333333, 0.157102, 0.005706288398682217, 0.005706288398682216, 3, 1, 0.4814053333333333, 1.444216, 0.0541864822020464, 0.05418648220204641, 0.02520991355911111, 0.226889222032, 0.0003092036947549378, 0.0003092036947549378),
    @("sCs", "Edil3", "Itgb3", "ECs", 3, 1, 1.092633, 3.277899, 0.119060464130005, 0.119060464130005, 3, 1, 3.778439, 11.335317, 0.4252971528324392, 0.4252971528324392, 4.128447139887, 37.156024258983, 0.05063607640939988, 0.05063607640939987),
    @("sCs", "Edil3", "Itgb3", "FAPs", 3, 1, 1.092633, 3.277899, 0.119060464130005, 0.119060464130005, 3, 1, 4.333403333333333, 13.00021, 0.4877633593505858, 0.4877633593505858, 4.734819484309999, 42.61337535879, 0.05807333194989114, 0.05807333194989114),
    @("sCs", "Edil3", "Itgb3", "M2", 3, 1, 1.092633, 3.277899, 0.119060464130005, 0.119060464130005, 3, 1, 0.2909853333333334, 0.8729560000000001, 0.03275300561492853, 0.03275300561492853, 0.317940177716, 2.861461599444, 0.00389958805016605, 0.00389958805016605),
    @("sCs", "Edil3", "Itgb3", "sCs", 3, 1, 1.092633, 3.277899, 0.119060464130005, 0.119060464130005, 3, 1, 0.4814053333333333, 1.444216, 0.0541864822020464, 0.05418648220204641, 0.5259993535759999, 4.733994182184, 0.006451467720547898, 0.006451467720547898)
)

$startRow = 2
for ($r = 0; $r -lt $rows.Length; $r++) {
    $values = $rows[$r]
    $excelRow = $startRow + $r
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value = $values[$c]
    }
}
